$d = $word.ActiveDocument

# The document's final paragraph ends with the "_GoBack" bookmark sitting
# right at the end of its text (after "...return to the destination. ").
# We want to add a brand new bullet/list paragraph after it containing
# "This solution meets all the goals. " and have the bookmark end up
# inside that new paragraph (at its end), exactly like the author typing
# the sentence and pressing Enter would produce.

$bm = $d.Bookmarks("_GoBack")
$bmRange = $bm.Range

# Insert the new sentence's text immediately before the bookmark. This
# keeps the bookmark glued to the end of the newly inserted text.
$bmRange.InsertBefore("This solution meets all the goals. ")

# Locate the start of the text we just inserted so we can split the
# paragraph right before it, turning it into its own (new) paragraph
# that inherits the same list formatting.
$lastPara = $d.Paragraphs.Last
$searchRange = $lastPara.Range.Duplicate
$searchRange.Find.Execute("This solution meets all the goals. ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$splitPos = $searchRange.Start

$splitRange = $d.Range($splitPos, $splitPos)
$splitRange.InsertParagraphBefore()
